$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart | Roof Tile
$ws.Cells.Item(19, 8).Value = 4151.121
$ws.Cells.Item(19, 9).Value = 640
$ws.Cells.Item(19, 10).Value = 5906.6816
$ws.Cells.Item(19, 11).Value = 640
$ws.Cells.Item(19, 12).Value = 5906.6816
$ws.Cells.Item(19, 13).Value = -465
$ws.Cells.Item(19, 14).Value = -6256.6816

# Row 64: Forged from the Void | Void Glue
$ws.Cells.Item(64, 8).Value = 11000
$ws.Cells.Item(64, 9).Value = 6666.6665
$ws.Cells.Item(64, 10).Value = 24000
$ws.Cells.Item(64, 11).Value = 6666.6665
$ws.Cells.Item(64, 12).Value = 24000
$ws.Cells.Item(64, 13).Value = -6418.6665
$ws.Cells.Item(64, 14).Value = -24496

# Row 67: Dodging the Draft (L) | Void Glue
$ws.Cells.Item(67, 8).Value = 11000
$ws.Cells.Item(67, 9).Value = 6666.6665
$ws.Cells.Item(67, 10).Value = 24000
$ws.Cells.Item(67, 11).Value = 6666.6665
$ws.Cells.Item(67, 12).Value = 24000
$ws.Cells.Item(67, 13).Value = -5808.6665
$ws.Cells.Item(67, 14).Value = -25716

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Cells.Item(45, 8).Value = 2053.7144
$ws.Cells.Item(45, 9).Value = 1479.5
$ws.Cells.Item(45, 11).Value = 1479.5
$ws.Cells.Item(45, 13).Value = -1102.5

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 2879.3572
$ws.Cells.Item(61, 9).Value = 1921.9362
$ws.Cells.Item(61, 11).Value = 1921.9362
$ws.Cells.Item(61, 13).Value = -1709.9362

# Row 63: Rivets Run through It | Mythrite Rivets
$ws.Cells.Item(63, 8).Value = 2999.5
$ws.Cells.Item(63, 9).Value = 2999.5
$ws.Cells.Item(63, 11).Value = 2999.5
$ws.Cells.Item(63, 13).Value = -2313.5

# Row 66: A Riveting Revival (L) | Mythrite Rivets
$ws.Cells.Item(66, 8).Value = 2999.5
$ws.Cells.Item(66, 9).Value = 2999.5
$ws.Cells.Item(66, 11).Value = 14997.5
$ws.Cells.Item(66, 13).Value = -11565.5

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 2879.3572
$ws.Cells.Item(136, 9).Value = 1921.9362
$ws.Cells.Item(136, 11).Value = 5765.8086
$ws.Cells.Item(136, 13).Value = -3215.8086

$ws = $wb.Worksheets.Item("BSM")
# Row 80: Unbreaker | Titanium Ingot
$ws.Cells.Item(80, 8).Value = 171.22223
$ws.Cells.Item(80, 9).Value = 145.8
$ws.Cells.Item(80, 10).Value = 181
$ws.Cells.Item(80, 11).Value = 145.8
$ws.Cells.Item(80, 12).Value = 181
$ws.Cells.Item(80, 13).Value = 852.2
$ws.Cells.Item(80, 14).Value = -2177

# Row 83: Attack on Titanium (L) | Titanium Ingot
$ws.Cells.Item(83, 8).Value = 171.22223
$ws.Cells.Item(83, 9).Value = 145.8
$ws.Cells.Item(83, 10).Value = 181
$ws.Cells.Item(83, 11).Value = 729
$ws.Cells.Item(83, 12).Value = 905
$ws.Cells.Item(83, 13).Value = 4263
$ws.Cells.Item(83, 14).Value = -10889

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Cells.Item(107, 8).Value = 3122.1562
$ws.Cells.Item(107, 9).Value = 3177
$ws.Cells.Item(107, 11).Value = 3177
$ws.Cells.Item(107, 13).Value = -1257

# Row 125: Archon of His Eye | High Durium Knives
$ws.Cells.Item(125, 8).Value = 199400
$ws.Cells.Item(125, 10).Value = 199400
$ws.Cells.Item(125, 12).Value = 199400
$ws.Cells.Item(125, 14).Value = -209240

# Row 141: Awl Dreams Come True | Ra'Kaznar Awl
$ws.Cells.Item(141, 8).Value = 53498
$ws.Cells.Item(141, 10).Value = 53498
$ws.Cells.Item(141, 12).Value = 53498
$ws.Cells.Item(141, 14).Value = -63858

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof | Ash Lumber
$ws.Cells.Item(16, 8).Value = 3422.6365
$ws.Cells.Item(16, 10).Value = 4409.1333
$ws.Cells.Item(16, 12).Value = 4409.1333
$ws.Cells.Item(16, 14).Value = -4983.1333

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Cells.Item(58, 8).Value = 2125.0513
$ws.Cells.Item(58, 9).Value = 1668.931
$ws.Cells.Item(58, 11).Value = 1668.931
$ws.Cells.Item(58, 13).Value = -1465.931

# Row 99: O Pine | Pine Lumber
$ws.Cells.Item(99, 8).Value = 2941.5386
$ws.Cells.Item(99, 9).Value = 2987.4092
$ws.Cells.Item(99, 11).Value = 2987.4092
$ws.Cells.Item(99, 13).Value = -1489.4092

# Row 113: Patient Patients | White Ash Lumber
$ws.Cells.Item(113, 8).Value = 3422.6365
$ws.Cells.Item(113, 10).Value = 4409.1333
$ws.Cells.Item(113, 12).Value = 4409.1333
$ws.Cells.Item(113, 14).Value = -8749.133300000001

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Cells.Item(126, 8).Value = 2941.5386
$ws.Cells.Item(126, 9).Value = 2987.4092
$ws.Cells.Item(126, 11).Value = 8962.2276
$ws.Cells.Item(126, 13).Value = -6492.2276

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Cells.Item(136, 8).Value = 2125.0513
$ws.Cells.Item(136, 9).Value = 1668.931
$ws.Cells.Item(136, 11).Value = 5006.793
$ws.Cells.Item(136, 13).Value = -2456.793

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food | Table Salt
$ws.Cells.Item(2, 8).Value = 521.2
$ws.Cells.Item(2, 9).Value = 647.55554
$ws.Cells.Item(2, 11).Value = 3885.33324
$ws.Cells.Item(2, 13).Value = -3772.33324

# Row 11: Putting the Squeeze On | Orange Juice
$ws.Cells.Item(11, 8).Value = 1001.7778
$ws.Cells.Item(11, 9).Value = 1838.6666
$ws.Cells.Item(11, 10).Value = 583.3333
$ws.Cells.Item(11, 11).Value = 5515.9998
$ws.Cells.Item(11, 12).Value = 1749.9999
$ws.Cells.Item(11, 13).Value = -5375.9998
$ws.Cells.Item(11, 14).Value = -2029.9999

# Row 12: Butter Me Up | Kukuru Butter
$ws.Cells.Item(12, 8).Value = 509.66666
$ws.Cells.Item(12, 9).Value = 429.33334
$ws.Cells.Item(12, 10).Value = 549.8333
$ws.Cells.Item(12, 11).Value = 1288.00002
$ws.Cells.Item(12, 12).Value = 1649.4999
$ws.Cells.Item(12, 13).Value = -1115.00002
$ws.Cells.Item(12, 14).Value = -1995.4999

# Row 26: A Grape Idea | Grape Juice
$ws.Cells.Item(26, 8).Value = 1626.4166
$ws.Cells.Item(26, 9).Value = 1879.7778
$ws.Cells.Item(26, 11).Value = 5639.3334
$ws.Cells.Item(26, 13).Value = -5351.3334

# Row 33: Cooking with Gas | Chicken Stock
$ws.Cells.Item(33, 8).Value = 928.5454999999999
$ws.Cells.Item(33, 9).Value = 33
$ws.Cells.Item(33, 10).Value = 3316.6667
$ws.Cells.Item(33, 11).Value = 198
$ws.Cells.Item(33, 12).Value = 19900.0002
$ws.Cells.Item(33, 13).Value = 85
$ws.Cells.Item(33, 14).Value = -20466.0002

# Row 122: Salt of the North | Northern Sea Salt
$ws.Cells.Item(122, 8).Value = 439.6
$ws.Cells.Item(122, 9).Value = 372
$ws.Cells.Item(122, 10).Value = 456.5
$ws.Cells.Item(122, 11).Value = 3348
$ws.Cells.Item(122, 12).Value = 4108.5
$ws.Cells.Item(122, 13).Value = -898
$ws.Cells.Item(122, 14).Value = -9008.5

$ws = $wb.Worksheets.Item("GSM")
# Row 11: A Ringing Success | Copper Ring
$ws.Cells.Item(11, 8).Value = 2185164.5
$ws.Cells.Item(11, 9).Value = 3635974.8
$ws.Cells.Item(11, 11).Value = 3635974.8
$ws.Cells.Item(11, 13).Value = -3635835.8

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Cells.Item(113, 8).Value = 250003360
$ws.Cells.Item(113, 9).Value = 250003360
$ws.Cells.Item(113, 11).Value = 250003360
$ws.Cells.Item(113, 13).Value = -250001190

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Cells.Item(122, 8).Value = 68385.92999999999
$ws.Cells.Item(122, 9).Value = 92544.82000000001
$ws.Cells.Item(122, 11).Value = 277634.46
$ws.Cells.Item(122, 13).Value = -275184.46

# Row 130: Planisphere to Paper | Chondrite Magitek Planisphere
$ws.Cells.Item(130, 8).Value = 24900
$ws.Cells.Item(130, 10).Value = 24900
$ws.Cells.Item(130, 12).Value = 24900
$ws.Cells.Item(130, 14).Value = -34940

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore | Hard Leather
$ws.Cells.Item(16, 8).Value = 1008.03925
$ws.Cells.Item(16, 9).Value = 929.0476
$ws.Cells.Item(16, 10).Value = 1376.6666
$ws.Cells.Item(16, 11).Value = 929.0476
$ws.Cells.Item(16, 12).Value = 1376.6666
$ws.Cells.Item(16, 13).Value = -759.0476
$ws.Cells.Item(16, 14).Value = -1716.6666

# Row 122: Hell on Leather | Gaja Leather
$ws.Cells.Item(122, 8).Value = 3191.5881
$ws.Cells.Item(122, 9).Value = 3152.4468
$ws.Cells.Item(122, 11).Value = 9457.340400000001
$ws.Cells.Item(122, 13).Value = -7007.340400000001

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Cells.Item(132, 8).Value = 3945.9167
$ws.Cells.Item(132, 9).Value = 2850.3572
$ws.Cells.Item(132, 10).Value = 5479.7
$ws.Cells.Item(132, 11).Value = 8551.071599999999
$ws.Cells.Item(132, 12).Value = 16439.1
$ws.Cells.Item(132, 13).Value = -6021.071599999999
$ws.Cells.Item(132, 14).Value = -21499.1

$ws = $wb.Worksheets.Item("WVR")
# Row 123: Helping Handwear | Fingerless Darkhempen Gloves of Healing
$ws.Cells.Item(123, 8).Value = 92713.5
$ws.Cells.Item(123, 10).Value = 92713.5
$ws.Cells.Item(123, 12).Value = 92713.5
$ws.Cells.Item(123, 14).Value = -102513.5

# Row 125: Color Coated | Almasty Serge Coat of Healing
$ws.Cells.Item(125, 8).Value = 57500
$ws.Cells.Item(125, 10).Value = 57500
$ws.Cells.Item(125, 12).Value = 57500
$ws.Cells.Item(125, 14).Value = -67340

# Row 131: A Better Bottom Line | AR-Caean Velvet Bottoms of Scouting
$ws.Cells.Item(131, 8).Value = 90000
$ws.Cells.Item(131, 10).Value = 90000
$ws.Cells.Item(131, 12).Value = 90000
$ws.Cells.Item(131, 14).Value = -100080

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 3278.3667
$ws.Cells.Item(132, 9).Value = 3319.6956
$ws.Cells.Item(132, 10).Value = 3142.5715
$ws.Cells.Item(132, 11).Value = 9959.086800000001
$ws.Cells.Item(132, 12).Value = 9427.7145
$ws.Cells.Item(132, 13).Value = -7429.086800000001
$ws.Cells.Item(132, 14).Value = -14487.7145
